# Update gh-pages to output generated at 456a3b4
# Applies the same numeric updates to both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 1207
    $ws.Range("G2").Value = 60

    $ws.Range("G5").Value = 60

    $ws.Range("F6").Value = 215

    $ws.Range("F10").Value = 5632

    $ws.Range("F11").Value = 4998
}
